{"js": "// Replace the date line and all 100 arithmetic-expression cells in the\n// table, in document order, with their new values (see commit diff).\n// Values below are listed in the same order the paragraphs occur in the\n// document body: index 0 is the title/date paragraph, indices 1..100 are\n// the table-cell paragraphs (row-major order through the 20-row x 5-column\n// table).\nconst oldValues = [\"2025-09-17 Wednesday\", \"53+2=\", \"49+42=\", \"86-25=\", \"40+41=\", \"82-49=\", \"13+74=\", \"29-27=\", \"35+49=\", \"50-23=\", \"44-8=\", \"9+57=\", \"72-42=\", \"25+7=\", \"69+11=\", \"23-15=\", \"15+39=\", \"3+5=\", \"60-42=\", \"85-14=\", \"91-29=\", \"39+23=\", \"24+38=\", \"2+53=\", \"52+47=\", \"18+59=\", \"92-27=\", \"12-7=\", \"45-29=\", \"71-39=\", \"58-46=\", \"10+50=\", \"32+62=\", \"54-17=\", \"24-17=\", \"46-17=\", \"93-56=\", \"37-4=\", \"82-9=\", \"65-4=\", \"37-34=\", \"5+21=\", \"73+21=\", \"47+46=\", \"45+46=\", \"34-25=\", \"67+4=\", \"36+33=\", \"76-10=\", \"14+50=\", \"17+51=\", \"98-62=\", \"74-30=\", \"32+7=\", \"32+46=\", \"66+13=\", \"34-14=\", \"51-33=\", \"15+22=\", \"42+2=\", \"44+8=\", \"39+28=\", \"5+1=\", \"20+47=\", \"31-15=\", \"82+1=\", \"34+38=\", \"72-2=\", \"22+67=\", \"74-26=\", \"92-45=\", \"95-19=\", \"70-23=\", \"92-0=\", \"70-12=\", \"47-13=\", \"61+32=\", \"64-43=\", \"59-18=\", \"60+37=\", \"59-44=\", \"85+1=\", \"42+33=\", \"54+13=\", \"22+31=\", \"60-42=\", \"11+17=\", \"3+93=\", \"1+14=\", \"21+59=\", \"69-42=\", \"29+3=\", \"72-2=\", \"8+22=\", \"15+14=\", \"20-11=\", \"89-24=\", \"92-63=\", \"28-23=\", \"23+71=\", \"29+37=\"];\nconst newValues = [\"2025-09-18 Thursday\", \"45-40=\", \"80+18=\", \"70+19=\", \"7-4=\", \"3+15=\", \"69-50=\", \"41-10=\", \"83-62=\", \"64-30=\", \"99-23=\", \"29+27=\", \"88-70=\", \"84-26=\", \"15+35=\", \"14-10=\", \"19+61=\", \"17+43=\", \"24+13=\", \"73-25=\", \"65-44=\", \"45-34=\", \"43+32=\", \"77-38=\", \"29+10=\", \"84-25=\", \"12+84=\", \"42+1=\", \"29+21=\", \"29+40=\", \"78-13=\", \"53+33=\", \"64+35=\", \"60+31=\", \"15+66=\", \"0+41=\", \"12+52=\", \"13+68=\", \"72-58=\", \"50-12=\", \"59-45=\", \"35-31=\", \"19+58=\", \"58-56=\", \"29+53=\", \"75+10=\", \"96-82=\", \"28+28=\", \"89-9=\", \"66-64=\", \"15+52=\", \"38-28=\", \"69-49=\", \"83-81=\", \"8+23=\", \"89-78=\", \"37+30=\", \"56-9=\", \"14+79=\", \"20-0=\", \"40-20=\", \"90-74=\", \"75-57=\", \"18+33=\", \"71-61=\", \"94-19=\", \"81-51=\", \"21-5=\", \"14+13=\", \"31-24=\", \"51+45=\", \"33+15=\", \"70+19=\", \"87-64=\", \"9+6=\", \"79-59=\", \"79-30=\", \"60-12=\", \"50-49=\", \"54+0=\", \"2+52=\", \"46+10=\", \"17+55=\", \"84+12=\", \"32+26=\", \"6+23=\", \"29+38=\", \"48+51=\", \"78-11=\", \"54-10=\", \"62+8=\", \"0+32=\", \"63-34=\", \"19+42=\", \"83-51=\", \"5+71=\", \"10+15=\", \"63-55=\", \"57+31=\", \"13-9=\", \"61-26=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} paragraphs, found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  // Defensive check: make sure we are about to replace the paragraph we\n  // expect (guards against unexpected document drift).\n  const current = para.text;\n  if (current !== oldValues[i]) {\n    throw new Error(\n      `Paragraph ${i} text mismatch: expected \"${oldValues[i]}\", found \"${current}\"`\n    );\n  }\n  if (newValues[i] !== oldValues[i]) {\n    para.insertText(newValues[i], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date/title line.\n$d = $word.ActiveDocument\n$titleOld = \"2025-09-17 Wednesday\"\n$titleNew = \"2025-09-18 Thursday\"\n$titlePara = $d.Paragraphs.Item(1)\n$titleCurrent = $titlePara.Range.Text.TrimEnd(\"`r\", \"`n\", [char]7)\nif ($titleCurrent -ne $titleOld) {\n    throw \"Title paragraph mismatch: expected '$titleOld', found '$titleCurrent'\"\n}\n$titlePara.Range.Text = $titleNew\n\n# Update every arithmetic-expression cell in the single 20x5 table, in\n# row-major (left-to-right, top-to-bottom) order -- the same order the\n# cells appear in the document.\n$oldValues = @(\n    \"53+2=\",\n    \"49+42=\",\n    \"86-25=\",\n    \"40+41=\",\n    \"82-49=\",\n    \"13+74=\",\n    \"29-27=\",\n    \"35+49=\",\n    \"50-23=\",\n    \"44-8=\",\n    \"9+57=\",\n    \"72-42=\",\n    \"25+7=\",\n    \"69+11=\",\n    \"23-15=\",\n    \"15+39=\",\n    \"3+5=\",\n    \"60-42=\",\n    \"85-14=\",\n    \"91-29=\",\n    \"39+23=\",\n    \"24+38=\",\n    \"2+53=\",\n    \"52+47=\",\n    \"18+59=\",\n    \"92-27=\",\n    \"12-7=\",\n    \"45-29=\",\n    \"71-39=\",\n    \"58-46=\",\n    \"10+50=\",\n    \"32+62=\",\n    \"54-17=\",\n    \"24-17=\",\n    \"46-17=\",\n    \"93-56=\",\n    \"37-4=\",\n    \"82-9=\",\n    \"65-4=\",\n    \"37-34=\",\n    \"5+21=\",\n    \"73+21=\",\n    \"47+46=\",\n    \"45+46=\",\n    \"34-25=\",\n    \"67+4=\",\n    \"36+33=\",\n    \"76-10=\",\n    \"14+50=\",\n    \"17+51=\",\n    \"98-62=\",\n    \"74-30=\",\n    \"32+7=\",\n    \"32+46=\",\n    \"66+13=\",\n    \"34-14=\",\n    \"51-33=\",\n    \"15+22=\",\n    \"42+2=\",\n    \"44+8=\",\n    \"39+28=\",\n    \"5+1=\",\n    \"20+47=\",\n    \"31-15=\",\n    \"82+1=\",\n    \"34+38=\",\n    \"72-2=\",\n    \"22+67=\",\n    \"74-26=\",\n    \"92-45=\",\n    \"95-19=\",\n    \"70-23=\",\n    \"92-0=\",\n    \"70-12=\",\n    \"47-13=\",\n    \"61+32=\",\n    \"64-43=\",\n    \"59-18=\",\n    \"60+37=\",\n    \"59-44=\",\n    \"85+1=\",\n    \"42+33=\",\n    \"54+13=\",\n    \"22+31=\",\n    \"60-42=\",\n    \"11+17=\",\n    \"3+93=\",\n    \"1+14=\",\n    \"21+59=\",\n    \"69-42=\",\n    \"29+3=\",\n    \"72-2=\",\n    \"8+22=\",\n    \"15+14=\",\n    \"20-11=\",\n    \"89-24=\",\n    \"92-63=\",\n    \"28-23=\",\n    \"23+71=\",\n    \"29+37=\"\n)\n$newValues = @(\n    \"45-40=\",\n    \"80+18=\",\n    \"70+19=\",\n    \"7-4=\",\n    \"3+15=\",\n    \"69-50=\",\n    \"41-10=\",\n    \"83-62=\",\n    \"64-30=\",\n    \"99-23=\",\n    \"29+27=\",\n    \"88-70=\",\n    \"84-26=\",\n    \"15+35=\",\n    \"14-10=\",\n    \"19+61=\",\n    \"17+43=\",\n    \"24+13=\",\n    \"73-25=\",\n    \"65-44=\",\n    \"45-34=\",\n    \"43+32=\",\n    \"77-38=\",\n    \"29+10=\",\n    \"84-25=\",\n    \"12+84=\",\n    \"42+1=\",\n    \"29+21=\",\n    \"29+40=\",\n    \"78-13=\",\n    \"53+33=\",\n    \"64+35=\",\n    \"60+31=\",\n    \"15+66=\",\n    \"0+41=\",\n    \"12+52=\",\n    \"13+68=\",\n    \"72-58=\",\n    \"50-12=\",\n    \"59-45=\",\n    \"35-31=\",\n    \"19+58=\",\n    \"58-56=\",\n    \"29+53=\",\n    \"75+10=\",\n    \"96-82=\",\n    \"28+28=\",\n    \"89-9=\",\n    \"66-64=\",\n    \"15+52=\",\n    \"38-28=\",\n    \"69-49=\",\n    \"83-81=\",\n    \"8+23=\",\n    \"89-78=\",\n    \"37+30=\",\n    \"56-9=\",\n    \"14+79=\",\n    \"20-0=\",\n    \"40-20=\",\n    \"90-74=\",\n    \"75-57=\",\n    \"18+33=\",\n    \"71-61=\",\n    \"94-19=\",\n    \"81-51=\",\n    \"21-5=\",\n    \"14+13=\",\n    \"31-24=\",\n    \"51+45=\",\n    \"33+15=\",\n    \"70+19=\",\n    \"87-64=\",\n    \"9+6=\",\n    \"79-59=\",\n    \"79-30=\",\n    \"60-12=\",\n    \"50-49=\",\n    \"54+0=\",\n    \"2+52=\",\n    \"46+10=\",\n    \"17+55=\",\n    \"84+12=\",\n    \"32+26=\",\n    \"6+23=\",\n    \"29+38=\",\n    \"48+51=\",\n    \"78-11=\",\n    \"54-10=\",\n    \"62+8=\",\n    \"0+32=\",\n    \"63-34=\",\n    \"19+42=\",\n    \"83-51=\",\n    \"5+71=\",\n    \"10+15=\",\n    \"63-55=\",\n    \"57+31=\",\n    \"13-9=\",\n    \"61-26=\"\n)\n\n$table = $d.Tables.Item(1)\n$rows = $table.Rows.Count\n$cols = $table.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $table.Cell($r, $c)\n        $current = $cell.Range.Text.TrimEnd(\"`r\", \"`n\", [char]7)\n        if ($current -ne $oldValues[$i]) {\n            throw \"Cell ($r,$c) mismatch: expected '$($oldValues[$i])', found '$current'\"\n        }\n        if ($newValues[$i] -ne $oldValues[$i]) {\n            $cell.Range.Text = $newValues[$i]\n        }\n        $i++\n    }\n}\n\n\"Updated $i cells\"\n"}
